# Updated cryptos list on Fri Sep 15 16:38:52 UTC 2023 with GitHub Actions
#
# Refreshes Price (D) / Volume(1h) (E) figures for the existing coins and
# inserts a new "BabyDogeCoin" entry at rank 46 (row 48), which pushes
# Cronos / EnergySwap / Mantle down one row each and drops USDD off the
# bottom of the A1:E51 table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text even when it "looks like" a number
# (e.g. "213.06", "0.0623") so Excel's auto-detect doesn't silently turn it
# into a numeric cell. NumberFormat is reset back to General and the style
# is reset to "Normal" afterwards so no stray formatting is left behind.
function Set-TextValue {
    param($Sheet, $Addr, $Val)
    $rng = $Sheet.Range($Addr)
    $rng.NumberFormat = "@"
    $rng.Value = $Val
    $rng.Style = "Normal"
}

# --- Price / Volume refresh for existing rows -------------------------
$ws.Range("D2").Value = "26.452.08"

$ws.Range("D3").Value = "1.626.49"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("E4").Value = "  +0.14%  "

Set-TextValue $ws "D5" "213.06"
$ws.Range("E5").Value = "  -0.07%  "

$ws.Range("E6").Value = "  +1.65%  "

$ws.Range("E7").Value = "  +0.12%  "

$ws.Range("E8").Value = "  -1.05%  "

Set-TextValue $ws "D9" "0.0623"
$ws.Range("E9").Value = "  +0.49%  "

Set-TextValue $ws "D10" "18.97"
$ws.Range("E10").Value = "  -0.62%  "

Set-TextValue $ws "D11" "0.0842"
$ws.Range("E11").Value = "  +0.94%  "

$ws.Range("D12").Value = "1.850.74"
$ws.Range("E12").Value = "  -0.71%  "

$ws.Range("D13").Value = "1.614.38"
$ws.Range("E13").Value = "  -1.26%  "

$ws.Range("E14").Value = "  +2.12%  "

Set-TextValue $ws "D15" "0.523"
$ws.Range("E15").Value = "  -0.11%  "

Set-TextValue $ws "D16" "64.40"
$ws.Range("E16").Value = "  +2.11%  "

$ws.Range("D17").Value = "26.442.66"
$ws.Range("E17").Value = "  -0.93%  "

$ws.Range("D18").Value = "0.0₃0739"
$ws.Range("E18").Value = "  +0.36%  "

Set-TextValue $ws "D19" "215.30"
$ws.Range("E19").Value = "  +2.41%  "

$ws.Range("E20").Value = "  +0.15%  "

Set-TextValue $ws "D21" "4.31"
$ws.Range("E21").Value = "  -0.23%  "

Set-TextValue $ws "D22" "6.24"
$ws.Range("E22").Value = "  +2.03%  "

Set-TextValue $ws "D23" "9.31"
$ws.Range("E23").Value = "  -0.77%  "

Set-TextValue $ws "D24" "1.99"
$ws.Range("E24").Value = "  +5.25%  "

Set-TextValue $ws "D25" "147.58"
$ws.Range("E25").Value = "  +0.68%  "

$ws.Range("E26").Value = "  +0.13%  "

$ws.Range("E27").Value = "  -0.52%  "

$ws.Range("E28").Value = "  +2.38%  "

$ws.Range("E29").Value = "  +0.89%  "

$ws.Range("E30").Value = "  -1.64%  "

$ws.Range("E31").Value = "  -1.35%  "

$ws.Range("E32").Value = "  +2.37%  "

Set-TextValue $ws "D33" "2.94"
$ws.Range("E33").Value = "  -0.45%  "

$ws.Range("E34").Value = "  -0.99%  "

$ws.Range("D35").Value = "1.219.26"
$ws.Range("E35").Value = "  +4.63%  "

$ws.Range("E36").Value = "  -1.38%  "

$ws.Range("E37").Value = "  +3.29%  "

$ws.Range("E38").Value = "  +0.12%  "

Set-TextValue $ws "D39" "0.797"
$ws.Range("E39").Value = "  -1.42%  "

Set-TextValue $ws "D40" "0.503"
$ws.Range("E40").Value = "  +0.13%  "

$ws.Range("E41").Value = "  -3.05%  "

Set-TextValue $ws "D42" "0.794"
$ws.Range("E42").Value = "  -0.13%  "

Set-TextValue $ws "D43" "5.37"
$ws.Range("E43").Value = "  +0.31%  "

$ws.Range("D44").Value = "1.761.35"
$ws.Range("E44").Value = "  -0.63%  "

Set-TextValue $ws "D45" "92.81"
$ws.Range("E45").Value = "  +0.24%  "

$ws.Range("E46").Value = "  +1.80%  "

Set-TextValue $ws "D47" "54.72"
$ws.Range("E47").Value = "  +0.12%  "

# --- Row 48-51 shuffle: new BabyDogeCoin row inserted, others shift down,
#     USDD (previously row 51) falls off the bottom of the table ---------
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0103"
$ws.Range("E48").Value = "  -0.66%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-TextValue $ws "D49" "0.0509"
$ws.Range("E49").Value = "  -0.59%  "

$ws.Range("B50").Value = "EnergySwap"
$ws.Range("C50").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D50" "7.54"
$ws.Range("E50").Value = "  -0.23%  "

$ws.Range("B51").Value = "Mantle"
$ws.Range("C51").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue $ws "D51" "0.407"
$ws.Range("E51").Value = "  -0.65%  "
